$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("A1").Value = 0.99915139503728889
$ws.Range("B1").Value = 0.00084860496271118668
$ws.Range("C1").Value = 0.00084860496271118668
$ws.Range("D1").Value = 0.00084860496271118733
$ws.Range("E1").Value = 0.00084860496271118733
$ws.Range("F1").Value = 0.00084860496271118506
$ws.Range("G1").Value = 0.00084860496271118441
$ws.Range("I1").Value = 0.00084860496271118733
$ws.Range("J1").Value = 0.00084860496271118668

# Row 2 updates
$ws.Range("A2").Value = 0.00084860496271118733
$ws.Range("D2").Value = 0.99915139503728889
$ws.Range("E2").Value = 0.99915139503728889
$ws.Range("H2").Value = 0.00084860496271118061
$ws.Range("J2").Value = 0.99915139503728889
